$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45203 -> 45204, i.e. 2023-10-04 -> 2023-10-05) for every data row (2..180).
$ws.Range("C2:C180").Value = 45204
